$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.9822660684585571
$ws.Range("B1").Value = 2.716087579727173
$ws.Range("C1").Value = 8.877663612365723
$ws.Range("D1").Value = 2.042148590087891
$ws.Range("E1").Value = 1.157161116600037
